$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 617
$ws.Range("A615:N615").Copy()
$ws.Range("A617:N617").PasteSpecial(-4122)
$ws.Cells.Item(617,1).Value = 45191.93974084491
$ws.Cells.Item(617,2).Value = "godis2002@naver.com"
$ws.Cells.Item(617,3).Value = "러시아학과"
$ws.Cells.Item(617,4).Value = 20211706
$ws.Cells.Item(617,5).Value = "김승겸"
$ws.Cells.Item(617,6).Value = "74:26"
$ws.Cells.Item(617,7).Value = 0.2
$ws.Cells.Item(617,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(617,9).Value = "952만 명"
$ws.Cells.Item(617,10).Value = 0.059
$ws.Cells.Item(617,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(617,12).Value = "Red"
$ws.Cells.Item(617,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A615:N615").Copy()
$ws.Range("A617:N617").PasteSpecial(-4122)

# Row 618
$ws.Range("A616:N616").Copy()
$ws.Range("A618:N618").PasteSpecial(-4122)
$ws.Cells.Item(618,1).Value = 45191.9549290625
$ws.Cells.Item(618,2).Value = "juh0611@naver.com"
$ws.Cells.Item(618,3).Value = "경영학과"
$ws.Cells.Item(618,4).Value = 20221051
$ws.Cells.Item(618,5).Value = "신주희"
$ws.Cells.Item(618,6).Value = "74:26"
$ws.Cells.Item(618,7).Value = 0.2
$ws.Cells.Item(618,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(618,9).Value = "952만 명"
$ws.Cells.Item(618,10).Value = 0.059
$ws.Cells.Item(618,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(618,12).Value = "Black"
$ws.Cells.Item(618,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A618:N618").PasteSpecial(-4122)

# Row 619
$ws.Range("A616:N616").Copy()
$ws.Range("A619:N619").PasteSpecial(-4122)
$ws.Cells.Item(619,1).Value = 45191.973035
$ws.Cells.Item(619,2).Value = "rlagkdud114@naver.com"
$ws.Cells.Item(619,3).Value = "정치행정학과"
$ws.Cells.Item(619,4).Value = 20221026
$ws.Cells.Item(619,5).Value = "김하영"
$ws.Cells.Item(619,6).Value = "77:23"
$ws.Cells.Item(619,7).Value = 0.2
$ws.Cells.Item(619,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(619,9).Value = "779만 명"
$ws.Cells.Item(619,10).Value = 0.151
$ws.Cells.Item(619,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(619,12).Value = "Black"
$ws.Cells.Item(619,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A619:N619").PasteSpecial(-4122)

# Row 620
$ws.Range("A616:N616").Copy()
$ws.Range("A620:N620").PasteSpecial(-4122)
$ws.Cells.Item(620,1).Value = 45191.98110314815
$ws.Cells.Item(620,2).Value = "sungjuwon1@gmail.com"
$ws.Cells.Item(620,3).Value = "컨텐츠IT전공"
$ws.Cells.Item(620,4).Value = 20215173
$ws.Cells.Item(620,5).Value = "성주원"
$ws.Cells.Item(620,6).Value = "74:26"
$ws.Cells.Item(620,7).Value = 0.2
$ws.Cells.Item(620,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(620,9).Value = "952만 명"
$ws.Cells.Item(620,10).Value = 0.059
$ws.Cells.Item(620,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(620,12).Value = "Black"
$ws.Cells.Item(620,14).Value = "모름/무응답"
$ws.Range("A616:N616").Copy()
$ws.Range("A620:N620").PasteSpecial(-4122)

# Row 621
$ws.Range("A615:N615").Copy()
$ws.Range("A621:N621").PasteSpecial(-4122)
$ws.Cells.Item(621,1).Value = 45191.98191899306
$ws.Cells.Item(621,2).Value = "janghangyeol0304@gmail.com"
$ws.Cells.Item(621,3).Value = "경제학과"
$ws.Cells.Item(621,4).Value = 20222838
$ws.Cells.Item(621,5).Value = "장한결"
$ws.Cells.Item(621,6).Value = "74:26"
$ws.Cells.Item(621,7).Value = 0.2
$ws.Cells.Item(621,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(621,9).Value = "952만 명"
$ws.Cells.Item(621,10).Value = 0.374
$ws.Cells.Item(621,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(621,12).Value = "Red"
$ws.Cells.Item(621,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A615:N615").Copy()
$ws.Range("A621:N621").PasteSpecial(-4122)

# Row 622
$ws.Range("A616:N616").Copy()
$ws.Range("A622:N622").PasteSpecial(-4122)
$ws.Cells.Item(622,1).Value = 45192.003076678244
$ws.Cells.Item(622,2).Value = "20233814@hallym.ac.kr"
$ws.Cells.Item(622,3).Value = "식품영양학과"
$ws.Cells.Item(622,4).Value = 20233814
$ws.Cells.Item(622,5).Value = "김정현"
$ws.Cells.Item(622,6).Value = "74:26"
$ws.Cells.Item(622,7).Value = 0.2
$ws.Cells.Item(622,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(622,9).Value = "952만 명"
$ws.Cells.Item(622,10).Value = 0.059
$ws.Cells.Item(622,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(622,12).Value = "Black"
$ws.Cells.Item(622,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A622:N622").PasteSpecial(-4122)

# Row 623
$ws.Range("A616:N616").Copy()
$ws.Range("A623:N623").PasteSpecial(-4122)
$ws.Cells.Item(623,1).Value = 45192.13073451389
$ws.Cells.Item(623,2).Value = "ridsigdog@gmail.com"
$ws.Cells.Item(623,3).Value = "경제학과"
$ws.Cells.Item(623,4).Value = 20202816
$ws.Cells.Item(623,5).Value = "박건민"
$ws.Cells.Item(623,6).Value = "74:26"
$ws.Cells.Item(623,7).Value = 0.25
$ws.Cells.Item(623,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(623,9).Value = "779만 명"
$ws.Cells.Item(623,10).Value = 0.151
$ws.Cells.Item(623,11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(623,12).Value = "Black"
$ws.Cells.Item(623,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A623:N623").PasteSpecial(-4122)

# Row 624
$ws.Range("A615:N615").Copy()
$ws.Range("A624:N624").PasteSpecial(-4122)
$ws.Cells.Item(624,1).Value = 45192.39833813657
$ws.Cells.Item(624,2).Value = "jimin4729@naver.com"
$ws.Cells.Item(624,3).Value = "법학과"
$ws.Cells.Item(624,4).Value = 20232747
$ws.Cells.Item(624,5).Value = "임지민"
$ws.Cells.Item(624,6).Value = "75:25"
$ws.Cells.Item(624,7).Value = 0.15
$ws.Cells.Item(624,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(624,9).Value = "166만 명"
$ws.Cells.Item(624,10).Value = 0.374
$ws.Cells.Item(624,11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(624,12).Value = "Red"
$ws.Cells.Item(624,13).Value = "모름/무응답"
$ws.Range("A615:N615").Copy()
$ws.Range("A624:N624").PasteSpecial(-4122)

# Row 625
$ws.Range("A615:N615").Copy()
$ws.Range("A625:N625").PasteSpecial(-4122)
$ws.Cells.Item(625,1).Value = 45192.49826432871
$ws.Cells.Item(625,2).Value = "jsy5233406@naver.com"
$ws.Cells.Item(625,3).Value = "식품영양학과"
$ws.Cells.Item(625,4).Value = 20233844
$ws.Cells.Item(625,5).Value = "정서연"
$ws.Cells.Item(625,6).Value = "75:25"
$ws.Cells.Item(625,7).Value = 0.15
$ws.Cells.Item(625,8).Value = "프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다."
$ws.Cells.Item(625,9).Value = "952만 명"
$ws.Cells.Item(625,10).Value = 0.059
$ws.Cells.Item(625,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(625,12).Value = "Red"
$ws.Cells.Item(625,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A615:N615").Copy()
$ws.Range("A625:N625").PasteSpecial(-4122)

# Row 626
$ws.Range("A616:N616").Copy()
$ws.Range("A626:N626").PasteSpecial(-4122)
$ws.Cells.Item(626,1).Value = 45192.5085250463
$ws.Cells.Item(626,2).Value = "040415kimdh@naver.com"
$ws.Cells.Item(626,3).Value = "소프트웨어학부"
$ws.Cells.Item(626,4).Value = 20235119
$ws.Cells.Item(626,5).Value = "김대현"
$ws.Cells.Item(626,6).Value = "74:26"
$ws.Cells.Item(626,7).Value = 0.2
$ws.Cells.Item(626,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(626,9).Value = "952만 명"
$ws.Cells.Item(626,10).Value = 0.059
$ws.Cells.Item(626,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(626,12).Value = "Black"
$ws.Cells.Item(626,14).Value = "모름/무응답"
$ws.Range("A616:N616").Copy()
$ws.Range("A626:N626").PasteSpecial(-4122)

# Row 627
$ws.Range("A615:N615").Copy()
$ws.Range("A627:N627").PasteSpecial(-4122)
$ws.Cells.Item(627,1).Value = 45192.53221329861
$ws.Cells.Item(627,2).Value = "20226760@hallym.ac.kr"
$ws.Cells.Item(627,3).Value = "Ai 의료융합"
$ws.Cells.Item(627,4).Value = 20226760
$ws.Cells.Item(627,5).Value = "이민홍"
$ws.Cells.Item(627,6).Value = "74:26"
$ws.Cells.Item(627,7).Value = 0.2
$ws.Cells.Item(627,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(627,9).Value = "952만 명"
$ws.Cells.Item(627,10).Value = 0.059
$ws.Cells.Item(627,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(627,12).Value = "Red"
$ws.Cells.Item(627,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A615:N615").Copy()
$ws.Range("A627:N627").PasteSpecial(-4122)

# Row 628
$ws.Range("A615:N615").Copy()
$ws.Range("A628:N628").PasteSpecial(-4122)
$ws.Cells.Item(628,1).Value = 45192.54229869213
$ws.Cells.Item(628,2).Value = "ans1929@gmail.com"
$ws.Cells.Item(628,3).Value = "법학과"
$ws.Cells.Item(628,4).Value = 20202750
$ws.Cells.Item(628,5).Value = "조유진"
$ws.Cells.Item(628,6).Value = "74:26"
$ws.Cells.Item(628,7).Value = 0.2
$ws.Cells.Item(628,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(628,9).Value = "952만 명"
$ws.Cells.Item(628,10).Value = 0.059
$ws.Cells.Item(628,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(628,12).Value = "Red"
$ws.Cells.Item(628,13).Value = "모름/무응답"
$ws.Range("A615:N615").Copy()
$ws.Range("A628:N628").PasteSpecial(-4122)

# Row 629
$ws.Range("A616:N616").Copy()
$ws.Range("A629:N629").PasteSpecial(-4122)
$ws.Cells.Item(629,1).Value = 45192.54367741898
$ws.Cells.Item(629,2).Value = "a01035025756@gmail.com"
$ws.Cells.Item(629,3).Value = "사회복지학부"
$ws.Cells.Item(629,4).Value = 20232342
$ws.Cells.Item(629,5).Value = "이건희"
$ws.Cells.Item(629,6).Value = "74:26"
$ws.Cells.Item(629,7).Value = 0.2
$ws.Cells.Item(629,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(629,9).Value = "779만 명"
$ws.Cells.Item(629,10).Value = 0.059
$ws.Cells.Item(629,11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(629,12).Value = "Black"
$ws.Cells.Item(629,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A629:N629").PasteSpecial(-4122)

# Row 630
$ws.Range("A615:N615").Copy()
$ws.Range("A630:N630").PasteSpecial(-4122)
$ws.Cells.Item(630,1).Value = 45192.56995640046
$ws.Cells.Item(630,2).Value = "kimhongik03@naver.com"
$ws.Cells.Item(630,3).Value = "심리학과"
$ws.Cells.Item(630,4).Value = 20232114
$ws.Cells.Item(630,5).Value = "김홍익"
$ws.Cells.Item(630,6).Value = "76:24"
$ws.Cells.Item(630,7).Value = 0.2
$ws.Cells.Item(630,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(630,9).Value = "779만 명"
$ws.Cells.Item(630,10).Value = 0.374
$ws.Cells.Item(630,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(630,12).Value = "Red"
$ws.Cells.Item(630,13).Value = "모름/무응답"
$ws.Range("A615:N615").Copy()
$ws.Range("A630:N630").PasteSpecial(-4122)

# Row 631
$ws.Range("A616:N616").Copy()
$ws.Range("A631:N631").PasteSpecial(-4122)
$ws.Cells.Item(631,1).Value = 45192.58355498842
$ws.Cells.Item(631,2).Value = "hyunjong9951@gmail.com"
$ws.Cells.Item(631,3).Value = "디스플레이"
$ws.Cells.Item(631,4).Value = 20183319
$ws.Cells.Item(631,5).Value = "이현종"
$ws.Cells.Item(631,6).Value = "74:26"
$ws.Cells.Item(631,7).Value = 0.2
$ws.Cells.Item(631,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(631,9).Value = "952만 명"
$ws.Cells.Item(631,10).Value = 0.059
$ws.Cells.Item(631,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(631,12).Value = "Black"
$ws.Cells.Item(631,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A631:N631").PasteSpecial(-4122)

# Row 632
$ws.Range("A615:N615").Copy()
$ws.Range("A632:N632").PasteSpecial(-4122)
$ws.Cells.Item(632,1).Value = 45192.59258246528
$ws.Cells.Item(632,2).Value = "serf0403@naver.com"
$ws.Cells.Item(632,3).Value = "바이오메디컬학과"
$ws.Cells.Item(632,4).Value = 20193646
$ws.Cells.Item(632,5).Value = "정예선"
$ws.Cells.Item(632,6).Value = "77:23"
$ws.Cells.Item(632,7).Value = 0.2
$ws.Cells.Item(632,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(632,9).Value = "166만 명"
$ws.Cells.Item(632,10).Value = 0.151
$ws.Cells.Item(632,11).Value = "상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다"
$ws.Cells.Item(632,12).Value = "Red"
$ws.Cells.Item(632,13).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A615:N615").Copy()
$ws.Range("A632:N632").PasteSpecial(-4122)

# Row 633
$ws.Range("A616:N616").Copy()
$ws.Range("A633:N633").PasteSpecial(-4122)
$ws.Cells.Item(633,1).Value = 45192.598411886574
$ws.Cells.Item(633,2).Value = "wnruddms@naver.com"
$ws.Cells.Item(633,3).Value = "식품영양학과"
$ws.Cells.Item(633,4).Value = 20203842
$ws.Cells.Item(633,5).Value = "주경은"
$ws.Cells.Item(633,6).Value = "74:26"
$ws.Cells.Item(633,7).Value = 0.2
$ws.Cells.Item(633,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(633,9).Value = "952만 명"
$ws.Cells.Item(633,10).Value = 0.059
$ws.Cells.Item(633,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(633,12).Value = "Black"
$ws.Cells.Item(633,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A633:N633").PasteSpecial(-4122)

# Row 634
$ws.Range("A615:N615").Copy()
$ws.Range("A634:N634").PasteSpecial(-4122)
$ws.Cells.Item(634,1).Value = 45192.61100075232
$ws.Cells.Item(634,2).Value = "freenix1001@naver.com"
$ws.Cells.Item(634,3).Value = "언어청각학부"
$ws.Cells.Item(634,4).Value = 20233903
$ws.Cells.Item(634,5).Value = "강서연"
$ws.Cells.Item(634,6).Value = "74:26"
$ws.Cells.Item(634,7).Value = 0.2
$ws.Cells.Item(634,8).Value = "OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다."
$ws.Cells.Item(634,9).Value = "779만 명"
$ws.Cells.Item(634,10).Value = 0.002
$ws.Cells.Item(634,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(634,12).Value = "Red"
$ws.Cells.Item(634,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A615:N615").Copy()
$ws.Range("A634:N634").PasteSpecial(-4122)

# Row 635
$ws.Range("A615:N615").Copy()
$ws.Range("A635:N635").PasteSpecial(-4122)
$ws.Cells.Item(635,1).Value = 45192.61873761574
$ws.Cells.Item(635,2).Value = "jonggwang0104@naver.com"
$ws.Cells.Item(635,3).Value = "식품영양학과"
$ws.Cells.Item(635,4).Value = 20183820
$ws.Cells.Item(635,5).Value = "박종광"
$ws.Cells.Item(635,6).Value = "74:26"
$ws.Cells.Item(635,7).Value = 0.2
$ws.Cells.Item(635,8).Value = "미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다."
$ws.Cells.Item(635,9).Value = "952만 명"
$ws.Cells.Item(635,10).Value = 0.059
$ws.Cells.Item(635,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(635,12).Value = "Red"
$ws.Cells.Item(635,13).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A615:N615").Copy()
$ws.Range("A635:N635").PasteSpecial(-4122)

# Row 636
$ws.Range("A616:N616").Copy()
$ws.Range("A636:N636").PasteSpecial(-4122)
$ws.Cells.Item(636,1).Value = 45192.62309835648
$ws.Cells.Item(636,2).Value = "1207dpwls@naver.com"
$ws.Cells.Item(636,3).Value = "법학과"
$ws.Cells.Item(636,4).Value = 20202751
$ws.Cells.Item(636,5).Value = "주예진"
$ws.Cells.Item(636,6).Value = "78:22"
$ws.Cells.Item(636,7).Value = 0.2
$ws.Cells.Item(636,8).Value = "우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다."
$ws.Cells.Item(636,9).Value = "38만 명"
$ws.Cells.Item(636,10).Value = 0.151
$ws.Cells.Item(636,11).Value = "중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다"
$ws.Cells.Item(636,12).Value = "Black"
$ws.Cells.Item(636,14).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A636:N636").PasteSpecial(-4122)

# Row 637
$ws.Range("A616:N616").Copy()
$ws.Range("A637:N637").PasteSpecial(-4122)
$ws.Cells.Item(637,1).Value = 45192.649936828704
$ws.Cells.Item(637,2).Value = "gwkang0330@gmail.com"
$ws.Cells.Item(637,3).Value = "소프트웨어학부"
$ws.Cells.Item(637,4).Value = 20235101
$ws.Cells.Item(637,5).Value = "강건우"
$ws.Cells.Item(637,6).Value = "77:23"
$ws.Cells.Item(637,7).Value = 0.15
$ws.Cells.Item(637,8).Value = "조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다"
$ws.Cells.Item(637,9).Value = "166만 명"
$ws.Cells.Item(637,10).Value = 0.002
$ws.Cells.Item(637,11).Value = "그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다"
$ws.Cells.Item(637,12).Value = "Black"
$ws.Cells.Item(637,14).Value = "국민부담률을 OECD 평균 수준으로 높여야 한다"
$ws.Range("A616:N616").Copy()
$ws.Range("A637:N637").PasteSpecial(-4122)

# Row 638
$ws.Range("A615:N615").Copy()
$ws.Range("A638:N638").PasteSpecial(-4122)
$ws.Cells.Item(638,1).Value = 45192.65257868056
$ws.Cells.Item(638,2).Value = "minjoo902@naver.com"
$ws.Cells.Item(638,3).Value = "금융재무학과"
$ws.Cells.Item(638,4).Value = 20192827
$ws.Cells.Item(638,5).Value = "김민주"
$ws.Cells.Item(638,6).Value = "74:26"
$ws.Cells.Item(638,7).Value = 0.2
$ws.Cells.Item(638,8).Value = "조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다"
$ws.Cells.Item(638,9).Value = "952만 명"
$ws.Cells.Item(638,10).Value = 0.059
$ws.Cells.Item(638,11).Value = "법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다."
$ws.Cells.Item(638,12).Value = "Red"
$ws.Cells.Item(638,13).Value = "국민부담률을 아일랜드 수준으로 낮춰야 한다"
$ws.Range("A615:N615").Copy()
$ws.Range("A638:N638").PasteSpecial(-4122)

